$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing account 002694089 / VITOR / 32592.53
# This is row 3 (row 1 = header "Conta/Nome/Saldo", row 2 = 004352384/BRASFORT,
# row 3 = 002694089/VITOR). Deleting the entire row shifts all following rows up.
$ws.Rows.Item(3).Delete()
